# Add semantic analysis for unused variable (#47)
#
# Populates the "Description" column (D) for the newly-added SemAnalyzer
# error/warning rows in the Errors table, and reclassifies row 26
# ("The variable is defined but not used") from an Error to a Warning
# (its FullId column is a formula that re-derives "W3014" automatically
# once the Severity cell changes).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Errors")

$ws.Range("D21").Value = "Invalid loop control statement context"
$ws.Range("D22").Value = "Invalid return context"
$ws.Range("D23").Value = "Not all path return a value"
$ws.Range("D24").Value = "The function requires a return value but not given"
$ws.Range("D25").Value = "The subprocedure cannot return a value but given"
$ws.Range("D26").Value = "The variable is defined but not used"

# Row 26 ("variable defined but not used") is a Warning, not an Error;
# F26's formula (LEFT(A26,1)&E26) recalculates E3014 -> W3014 on its own.
$ws.Range("A26").Value = "Warning"

# Move the active selection to where the edits were made.
$null = $ws.Activate()
$null = $ws.Range("D27").Select()
